$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 9852.227999999999
$ws.Cells.Item(19, 10).Value = 13262.4375
$ws.Cells.Item(19, 12).Value = 13262.4375
$ws.Cells.Item(19, 14).Value = -13612.4375

$ws.Cells.Item(112, 8).Value = 5129.4736
$ws.Cells.Item(112, 9).Value = 400
$ws.Cells.Item(112, 10).Value = 5392.222
$ws.Cells.Item(112, 11).Value = 1200
$ws.Cells.Item(112, 12).Value = 16176.666
$ws.Cells.Item(112, 13).Value = -92
$ws.Cells.Item(112, 14).Value = -18392.666

$ws.Cells.Item(115, 8).Value = 2250.4736
$ws.Cells.Item(115, 9).Value = 1417.6666
$ws.Cells.Item(115, 11).Value = 4252.9998
$ws.Cells.Item(115, 13).Value = -2685.9998

$ws.Cells.Item(116, 8).Value = 35335
$ws.Cells.Item(116, 9).Value = 35335
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 35335
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).Value = -31893
$ws.Cells.Item(116, 14).ClearContents()

$ws.Cells.Item(129, 8).Value = 1166.4032
$ws.Cells.Item(129, 9).Value = 557
$ws.Cells.Item(129, 10).Value = 1312.66
$ws.Cells.Item(129, 11).Value = 1671
$ws.Cells.Item(129, 12).Value = 3937.98
$ws.Cells.Item(129, 13).Value = 3329
$ws.Cells.Item(129, 14).Value = -13937.98

$ws.Cells.Item(132, 8).Value = 2011.4073
$ws.Cells.Item(132, 9).Value = 1727.102
$ws.Cells.Item(132, 11).Value = 5181.306
$ws.Cells.Item(132, 13).Value = -2651.306

$ws.Cells.Item(137, 8).Value = 2559.0154
$ws.Cells.Item(137, 9).Value = 2342.0408
$ws.Cells.Item(137, 10).Value = 3223.5
$ws.Cells.Item(137, 11).Value = 7026.1224
$ws.Cells.Item(137, 12).Value = 9670.5
$ws.Cells.Item(137, 13).Value = -4476.1224
$ws.Cells.Item(137, 14).Value = -14770.5

$ws.Cells.Item(138, 8).Value = 2568.1
$ws.Cells.Item(138, 9).Value = 1247.8636
$ws.Cells.Item(138, 10).Value = 3830.9348
$ws.Cells.Item(138, 11).Value = 3743.5908
$ws.Cells.Item(138, 12).Value = 11492.8044
$ws.Cells.Item(138, 13).Value = 1396.4092
$ws.Cells.Item(138, 14).Value = -21772.8044

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 22730166
$ws.Cells.Item(61, 9).Value = 43480590
$ws.Cells.Item(61, 11).Value = 43480590
$ws.Cells.Item(61, 13).Value = -43480378

$ws.Cells.Item(136, 8).Value = 22730166
$ws.Cells.Item(136, 9).Value = 43480590
$ws.Cells.Item(136, 11).Value = 130441770
$ws.Cells.Item(136, 13).Value = -130439220

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(100, 8).Value = 85000
$ws.Cells.Item(100, 10).Value = 85000
$ws.Cells.Item(100, 12).Value = 85000
$ws.Cells.Item(100, 14).Value = -87164

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 48000.5
$ws.Cells.Item(4, 9).Value = 25000
$ws.Cells.Item(4, 11).Value = 25000
$ws.Cells.Item(4, 13).Value = -24888

$ws.Cells.Item(16, 8).Value = 633.1667
$ws.Cells.Item(16, 9).Value = 633.1667
$ws.Cells.Item(16, 11).Value = 633.1667
$ws.Cells.Item(16, 13).Value = -346.1667

$ws.Cells.Item(31, 8).Value = 7966.8716
$ws.Cells.Item(31, 9).Value = 1385.5294
$ws.Cells.Item(31, 10).Value = 13052.454
$ws.Cells.Item(31, 11).Value = 1385.5294
$ws.Cells.Item(31, 12).Value = 13052.454
$ws.Cells.Item(31, 13).Value = -1090.5294
$ws.Cells.Item(31, 14).Value = -13642.454

$ws.Cells.Item(34, 8).Value = 7966.8716
$ws.Cells.Item(34, 9).Value = 1385.5294
$ws.Cells.Item(34, 10).Value = 13052.454
$ws.Cells.Item(34, 11).Value = 1385.5294
$ws.Cells.Item(34, 12).Value = 13052.454
$ws.Cells.Item(34, 13).Value = -1183.5294
$ws.Cells.Item(34, 14).Value = -13456.454

$ws.Cells.Item(45, 8).Value = 6309
$ws.Cells.Item(45, 9).Value = 6309
$ws.Cells.Item(45, 11).Value = 6309
$ws.Cells.Item(45, 13).Value = -5716

$ws.Cells.Item(94, 8).Value = 166668180
$ws.Cells.Item(94, 9).Value = 1000000000
$ws.Cells.Item(94, 10).Value = 1807.4
$ws.Cells.Item(94, 11).Value = 1000000000
$ws.Cells.Item(94, 12).Value = 1807.4
$ws.Cells.Item(94, 13).Value = -999999549
$ws.Cells.Item(94, 14).Value = -2709.4

$ws.Cells.Item(99, 8).Value = 2000
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 2000
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = 2000
$ws.Cells.Item(99, 13).ClearContents()
$ws.Cells.Item(99, 14).Value = -4996

$ws.Cells.Item(113, 8).Value = 633.1667
$ws.Cells.Item(113, 9).Value = 633.1667
$ws.Cells.Item(113, 11).Value = 633.1667
$ws.Cells.Item(113, 13).Value = 1536.8333

$ws.Cells.Item(126, 8).Value = 2000
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 2000
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 6000
$ws.Cells.Item(126, 13).Value = -10940
$ws.Cells.Item(126, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 21213238
$ws.Cells.Item(132, 9).Value = 21740016
$ws.Cells.Item(132, 10).Value = 18520814
$ws.Cells.Item(132, 11).Value = 65220048
$ws.Cells.Item(132, 12).Value = 55562442
$ws.Cells.Item(132, 13).Value = -65217518
$ws.Cells.Item(132, 14).Value = -55567502

$ws.Cells.Item(134, 8).Value = 3785.9565
$ws.Cells.Item(134, 9).Value = 3706.6052
$ws.Cells.Item(134, 10).Value = 4162.875
$ws.Cells.Item(134, 11).Value = 11119.8156
$ws.Cells.Item(134, 12).Value = 12488.625
$ws.Cells.Item(134, 13).Value = -8584.8156
$ws.Cells.Item(134, 14).Value = -17558.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 21819016
$ws.Cells.Item(4, 9).Value = 15556352
$ws.Cells.Item(4, 10).Value = 50001000
$ws.Cells.Item(4, 11).Value = 46669056
$ws.Cells.Item(4, 12).Value = 150003000
$ws.Cells.Item(4, 13).Value = -46668944
$ws.Cells.Item(4, 14).Value = -150003224

$ws.Cells.Item(44, 8).Value = 261.44446
$ws.Cells.Item(44, 9).Value = 261.44446
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 11).Value = 784.33338
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 13).Value = -386.33338
$ws.Cells.Item(44, 14).ClearContents()

$ws.Cells.Item(113, 8).Value = 1143.4706
$ws.Cells.Item(113, 9).Value = 628.3333
$ws.Cells.Item(113, 10).Value = 1424.4546
$ws.Cells.Item(113, 11).Value = 1884.9999
$ws.Cells.Item(113, 12).Value = 4273.3638
$ws.Cells.Item(113, 13).Value = 285.0001
$ws.Cells.Item(113, 14).Value = -8613.363799999999

$ws.Cells.Item(122, 8).Value = 2921.7673
$ws.Cells.Item(122, 10).Value = 6498.5
$ws.Cells.Item(122, 12).Value = 58486.5
$ws.Cells.Item(122, 14).Value = -63386.5

$ws.Cells.Item(137, 8).Value = 9812830
$ws.Cells.Item(137, 9).Value = 41692170
$ws.Cells.Item(137, 10).Value = 3802.5386
$ws.Cells.Item(137, 11).Value = 125076510
$ws.Cells.Item(137, 12).Value = 11407.6158
$ws.Cells.Item(137, 13).Value = -125071410
$ws.Cells.Item(137, 14).Value = -21607.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 3000

$ws.Cells.Item(74, 8).Value = 43066.555
$ws.Cells.Item(74, 10).Value = 43066.555
$ws.Cells.Item(74, 12).Value = 43066.555
$ws.Cells.Item(74, 14).Value = -44938.555

$ws.Cells.Item(75, 8).Value = 26495
$ws.Cells.Item(75, 10).Value = 26495
$ws.Cells.Item(75, 12).Value = 26495
$ws.Cells.Item(75, 13).Value = -28243

$ws.Cells.Item(77, 8).Value = 43066.555
$ws.Cells.Item(77, 10).Value = 43066.555
$ws.Cells.Item(77, 12).Value = 129199.665
$ws.Cells.Item(77, 14).Value = -138559.665

$ws.Cells.Item(78, 8).Value = 26495
$ws.Cells.Item(78, 10).Value = 26495
$ws.Cells.Item(78, 12).Value = 79485
$ws.Cells.Item(78, 14).Value = -88221

$ws.Cells.Item(126, 8).Value = 2615.1428
$ws.Cells.Item(126, 9).Value = 2615.1428
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 7845.428400000001
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -5375.428400000001
$ws.Cells.Item(126, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 3757.879
$ws.Cells.Item(132, 9).Value = 3635.7273
$ws.Cells.Item(132, 10).Value = 4002.182
$ws.Cells.Item(132, 11).Value = 10907.1819
$ws.Cells.Item(132, 12).Value = 12006.546
$ws.Cells.Item(132, 13).Value = -8377.1819
$ws.Cells.Item(132, 14).Value = -17066.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 4046668.8
$ws.Cells.Item(2, 10).Value = 4046668.8
$ws.Cells.Item(2, 12).Value = 4046668.8
$ws.Cells.Item(2, 14).Value = -4046892.8

$ws.Cells.Item(132, 8).Value = 2287.587
$ws.Cells.Item(132, 10).Value = 2781.9375
$ws.Cells.Item(132, 12).Value = 8345.8125
$ws.Cells.Item(132, 14).Value = -13405.8125

$ws.Cells.Item(133, 8).Value = 38550
$ws.Cells.Item(133, 10).Value = 38550
$ws.Cells.Item(133, 12).Value = 38550
$ws.Cells.Item(133, 14).Value = -43610

$ws.Cells.Item(136, 8).Value = 2605775.8
$ws.Cells.Item(136, 9).Value = 1136.5518
$ws.Cells.Item(136, 10).Value = 4763905
$ws.Cells.Item(136, 11).Value = 3409.6554
$ws.Cells.Item(136, 12).Value = 14291715
$ws.Cells.Item(136, 13).Value = -859.6553999999996
$ws.Cells.Item(136, 14).Value = -14296815

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 102073000
$ws.Cells.Item(2, 9).Value = 6833333.5
$ws.Cells.Item(2, 10).Value = 142890000
$ws.Cells.Item(2, 11).Value = 6833333.5
$ws.Cells.Item(2, 12).Value = 142890000
$ws.Cells.Item(2, 13).Value = -6833221.5
$ws.Cells.Item(2, 14).Value = -142890224

$ws.Cells.Item(136, 8).Value = 1802.1
$ws.Cells.Item(136, 10).Value = 3329.2307
$ws.Cells.Item(136, 12).Value = 9987.6921
$ws.Cells.Item(136, 14).Value = -15087.6921
